# 9th Stab - Cosmetic Changes
# Insert two new "week" columns (Jun_17, Jun_15) to the left of the existing
# Jun_13 / Jun_10 columns, shifting the old B (Jun_13) and C (Jun_10) data
# over to D and E respectively, and filling the two new columns with the
# default "UN" marker used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at B:C -- this shifts existing column B -> D and
# existing column C -> E automatically, carrying their values/formatting.
$ws.Range("B:C").EntireColumn.Insert()

# New header row (row 1) values for the freshly inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the rest of the new B & C columns (rows 2-27) with the same default
# placeholder value ("UN") used elsewhere on the sheet.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Restore/apply the column width formatting (8 characters wide, custom width)
# across the now-adjacent C, D, E columns.
$ws.Range("C:E").ColumnWidth = 7.1666666666666667
